$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 372, pushing the existing rows 372-407 down
# to 373-408 (this also extends dimension to A1:R408 automatically).
$ws.Rows.Item(372).Insert()

# Populate the newly inserted row 372 with the new weekly record.
$ws.Cells.Item(372, 1).Value = 10
$ws.Cells.Item(372, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(372, 3).Value = "La Araucanía"
$ws.Cells.Item(372, 4).Value = 45132
$ws.Cells.Item(372, 5).Value = 9
$ws.Cells.Item(372, 6).Value = 100112039
$ws.Cells.Item(372, 7).Value = "Ciboulette"
$ws.Cells.Item(372, 8).Value = "Sin especificar"
$ws.Cells.Item(372, 9).Value = "Primera"
$ws.Cells.Item(372, 10).Value = 50
$ws.Cells.Item(372, 11).Value = 7000
$ws.Cells.Item(372, 12).Value = 7000
$ws.Cells.Item(372, 13).Value = 7000
$ws.Cells.Item(372, 14).Value = "$/docena de atados"
$ws.Cells.Item(372, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(372, 16).Value = 2333
$ws.Cells.Item(372, 17).Value = 3
$ws.Cells.Item(372, 18).Value = "Hortaliza"
